$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Occurrence 1 (table cell): "DEEP LEARNING FOR BEGGINERS: NEURAL
# NETWORKS IN R STUDIO." -> fix typo BEGGINERS -> BEGINNERS, then
# split the run into 5 pieces matching the target edit.
# ---------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("DEEP LEARNING FOR BEGGINERS: NEURAL NETWORKS IN R STUDIO.", $true, $false, $false, $false, $false, $true, 1, $false, "DEEP LEARNING FOR BEGINNERS: NEURAL NETWORKS IN R STUDIO.", 1)

$r1b = $d.Content
$r1b.Find.Execute("DEEP LEARNING FOR BEGINNERS: NEURAL NETWORKS IN R STUDIO.")
$s1 = $r1b.Start

$p1a = $d.Range($s1+0, $s1+13)
$p1a.Font.Size = 13
$p1a.Font.Size = 12

$p1b = $d.Range($s1+13, $s1+21)
$p1b.Font.Size = 13
$p1b.Font.Size = 12

$p1c = $d.Range($s1+21, $s1+23)
$p1c.Font.Size = 13
$p1c.Font.Size = 12

$p1d = $d.Range($s1+23, $s1+24)
$p1d.Font.Size = 13
$p1d.Font.Size = 12

# ---------------------------------------------------------------
# Occurrence 2 (paragraph): "CERTIFICATION COURSE NAME:- DEEP
# LEARNING FOR BEGGINERS: NEURAL NETWORKS IN R STUDIO." -> same
# typo fix, split into 7 pieces matching the target edit.
# ---------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("CERTIFICATION COURSE NAME:- DEEP LEARNING FOR BEGGINERS: NEURAL NETWORKS IN R STUDIO.", $true, $false, $false, $false, $false, $true, 1, $false, "CERTIFICATION COURSE NAME:- DEEP LEARNING FOR BEGINNERS: NEURAL NETWORKS IN R STUDIO.", 1)

$r2b = $d.Content
$r2b.Find.Execute("CERTIFICATION COURSE NAME:- DEEP LEARNING FOR BEGINNERS: NEURAL NETWORKS IN R STUDIO.")
$s2 = $r2b.Start

$p2a = $d.Range($s2+0, $s2+18)
$p2a.Font.Size = 13
$p2a.Font.Size = 12

$p2b = $d.Range($s2+18, $s2+25)
$p2b.Font.Size = 13
$p2b.Font.Size = 12

$p2c = $d.Range($s2+25, $s2+27)
$p2c.Font.Size = 13
$p2c.Font.Size = 12

$p2d = $d.Range($s2+27, $s2+49)
$p2d.Font.Size = 13
$p2d.Font.Size = 12

$p2e = $d.Range($s2+49, $s2+50)
$p2e.Font.Size = 13
$p2e.Font.Size = 12

$p2f = $d.Range($s2+50, $s2+51)
$p2f.Font.Size = 13
$p2f.Font.Size = 12
